$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.529.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "'2.602.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'539.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.28%  "
$ws.Range("D6").Value = "'141.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("D13").Value = "'3.062.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'59.462.18"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").Value = "'20.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "'2.641.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "'341.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "'4.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.47%  "
$ws.Range("D20").Value = "'10.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  -1.86%  "
$ws.Range("E22").Value = "  -0.07%  "
$ws.Range("D23").Value = "'67.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.30%  "
$ws.Range("D24").Value = "'0.408"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.25%  "
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Value = "'7.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.10%  "
$ws.Range("D28").Value = "'0.0₃0743"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("E30").Value = "  +6.01%  "
$ws.Range("D31").Value = "'5.83"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").Value = "'18.78"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Value = "'150.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.45%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").Value = "'0.845"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.16%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").Value = "'0.824"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.01%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "'272.48"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  -0.40%  "
$ws.Range("D43").Value = "'10.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "'0.0950"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").Value = "'18.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.17%  "
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").Value = "'1.937.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("D50").Value = "'111.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("E51").Value = "  +2.16%  "
